# Calculated ROC AUC with probabilities
# Flip the predicted class label ("no"/"yes") in column A for the rows
# below that were recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row number -> new text value for column A
$changes = @{
    5   = "yes"
    21  = "no"
    48  = "yes"
    309 = "yes"
    319 = "yes"
    348 = "yes"
    403 = "yes"
    428 = "yes"
    488 = "yes"
    539 = "yes"
    542 = "no"
    648 = "no"
    659 = "yes"
    735 = "no"
    741 = "no"
    758 = "no"
    762 = "yes"
    773 = "yes"
    784 = "no"
    804 = "yes"
    808 = "no"
    825 = "no"
    870 = "yes"
    873 = "yes"
    893 = "yes"
}

foreach ($row in $changes.Keys) {
    $ws.Range("A$row").Value = $changes[$row]
}
